$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Bio" column (column S) with header + value
$ws.Range("S1").Value = "Bio"
$ws.Range("S2").Value = "As a versatile Digital Content Creator, I bring ideas to life across the digital landscape, specializing in end-to-end media production, from compelling Design and professional-grade Video Editing to pioneering the use of AI content generation. I thrive on transforming complex concepts into captivating, high-impact digital experiences that drive engagement and tell unforgettable stories."

# The previously auto-fit "Brand logos" column (R) is narrowed to a fixed custom width
$ws.Columns.Item(18).ColumnWidth = 22

# Update view state to match where the author left the selection/scroll
$ws.Range("S15").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 11
